# Commit: "Restored from revision of admin on 12/14/2020 09:53:48 AM.TEST
# Author: admin. Type: SAVE." -- the only functional change in the target
# revision is the value stored in C10 on the "Rules" sheet, which goes
# from 18 to 1. (Everything else in the supplied diff -- numFmts/xfId
# additions, the customWidth="true"/collapsed="false" col attribute
# reformatting, the t="n"/1.0 float spelling -- is re-serialization noise
# from the tool that produced the XML snapshot, not a distinct edit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
